$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for Wins / Losses / Ties in columns AD, AE, AF (row 1)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the header style used by the rest of row 1 (e.g. AC1) by copying its formatting
$headerStyleSource = $ws.Range("AC1")
$headerStyleSource.Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the team record (Wins, Losses, Ties) for every data row
$lastRow = 51
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 64
    $ws.Cells.Item($r, 31).Value = 98
    $ws.Cells.Item($r, 32).Value = 0
}
